# LMS_Request_Details.xlsx update:
#  - add Sheet2..Sheet5 (after Sheet1, in order)
#  - Sheet1: add column B ("invalidProgramId" / "209876" as text, matching
#    the existing quote-prefixed-text style used for "345678" in A2)
#  - Sheet1: widen column B, move the selection to C5
#  - Sheet4: leave a lingering selection at E1
#  - return focus to Sheet1 so it stays the active tab

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item(1)

# Insert four new sheets right after Sheet1, keeping them in numeric order.
$wb.Worksheets.Add($null, $sheet1) | Out-Null
$sheet2 = $wb.Worksheets.Item(2)
$wb.Worksheets.Add($null, $sheet2) | Out-Null
$sheet3 = $wb.Worksheets.Item(3)
$wb.Worksheets.Add($null, $sheet3) | Out-Null
$sheet4 = $wb.Worksheets.Item(4)
$wb.Worksheets.Add($null, $sheet4) | Out-Null

# Visit Sheet4 and leave its selection at E1.
$sheet4 = $wb.Worksheets.Item("Sheet4")
$sheet4.Activate()
$sheet4.Range("E1").Select() | Out-Null

# Back to Sheet1: add the new "invalidProgramId" / "209876" column.
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Activate()

$sheet1.Range("B1").Value = "invalidProgramId"
$sheet1.Range("B2").Value = "'209876"

$sheet1.Columns.Item(2).ColumnWidth = 17.166666666666668

$sheet1.Range("C5").Select() | Out-Null
